# Apply the "added 3 new FX" update to the Audio Asset List workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: Woosh.wav (new ambient/FX entry, now "Done") ---
$ws.Range("A3").Value = "Woosh.wav"
$ws.Range("D3").Value = "Done"
$ws.Range("E3").Value = "I held the mic up to my desk fan"

# --- Row 4: big-explosion -> .wav, now "Done", new note ---
$ws.Range("A4").Value = "big-explosion.wav"
$ws.Range("D4").Value = "Done"
$ws.Range("E4").Value = "replace with its own heavier reverb version"

# --- Row 5: Proceed.wav (new interface entry, now "Done") ---
$ws.Range("A5").Value = "Proceed.wav"
$ws.Range("D5").Value = "Done"

# --- Row 6: Skip.wav (new interface entry, now "Done") ---
$ws.Range("A6").Value = "Skip.wav"
$ws.Range("D6").Value = "Done"

# --- Row 7: Pew.mp3 -> Pew.wav, now "Done" ---
$ws.Range("A7").Value = "Pew.wav"
$ws.Range("D7").Value = "Done"

# --- Row 8: big-explosion.mp3 -> big-explosion.wav, now "Done" ---
$ws.Range("A8").Value = "big-explosion.wav"
$ws.Range("D8").Value = "Done"

# --- Row 9: small-explosion.mp3 -> small-explosion.wav, now "Done" ---
$ws.Range("A9").Value = "small-explosion.wav"
$ws.Range("D9").Value = "Done"

# --- Row 12: unchanged text; only its shared-string index shifts ---

# --- Row 15: remove the "Not Done" status (blank cell) ---
$ws.Range("D15").Value = ""

# --- Rows 16,18,20,22,26,29,32,34: status -> "In progress" ---
$ws.Range("D16").Value = "In progress"
$ws.Range("D18").Value = "In progress"
$ws.Range("D20").Value = "In progress"
$ws.Range("D22").Value = "In progress"
$ws.Range("D26").Value = "In progress"
$ws.Range("D29").Value = "In progress"
$ws.Range("D32").Value = "In progress"
$ws.Range("D34").Value = "In progress"

# --- update selection to match the new active cell in the diff ---
$ws.Range("E7").Select()
